# Applies the FortiGate base-config tidy-up edit:
#  - Row 2 (HQ): rotate the fortios_access_token value, drop the cell's
#    wrap-text formatting
#  - Row 2 (HQ): rename the default outbound policy + its source interface
#  - Row 3 (Spoke1): clear the fortios_access_token cell (C3)
#  - Selection / view housekeeping to match the saved state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: rotate the HQ access token and clear the wrap-text formatting ---
$ws.Range("C2").Value = "GN1cH1bmb4t6yGt3qb7fwmknmwmNff"
$ws.Range("C2").WrapText = $false

# --- Row 2: rename default outbound policy / source interface ---
$ws.Range("W2").Value = "Default-Outbound"
$ws.Range("X2").Value = "LAN"

# --- Row 3: clear the Spoke1 access-token cell ---
$ws.Range("C3").ClearContents()

# --- View / selection housekeeping ---
$ws.Range("C2").Select() | Out-Null
